# Weekly update: insert two new price observations (rows 260-261) for
# Coliflor at Terminal Hortofrutícola Agro Chillán, shifting the rest of
# the table down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 260-261; existing rows 260.. shift down to 262..
$ws.Rows("260:261").Insert()

# --- New row 260 ---
$ws.Cells.Item(260, 1).Value = 7
$ws.Cells.Item(260, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(260, 3).Value = "Ñuble"
$ws.Cells.Item(260, 4).Value = 44876
$ws.Cells.Item(260, 5).Value = 16
$ws.Cells.Item(260, 6).Value = 100112008
$ws.Cells.Item(260, 7).Value = "Coliflor"
$ws.Cells.Item(260, 8).Value = "Sin especificar"
$ws.Cells.Item(260, 9).Value = "Primera"
$ws.Cells.Item(260, 10).Value = 300
$ws.Cells.Item(260, 11).Value = 700
$ws.Cells.Item(260, 12).Value = 800
$ws.Cells.Item(260, 13).Value = 750
$ws.Cells.Item(260, 14).Value = "`$/unidad"
$ws.Cells.Item(260, 15).Value = "Región del Maule"
$ws.Cells.Item(260, 16).Value = 750
$ws.Cells.Item(260, 17).Value = 1
$ws.Cells.Item(260, 18).Value = "Hortaliza"

# --- New row 261 ---
$ws.Cells.Item(261, 1).Value = 7
$ws.Cells.Item(261, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(261, 3).Value = "Ñuble"
$ws.Cells.Item(261, 4).Value = 44876
$ws.Cells.Item(261, 5).Value = 16
$ws.Cells.Item(261, 6).Value = 100112008
$ws.Cells.Item(261, 7).Value = "Coliflor"
$ws.Cells.Item(261, 8).Value = "Sin especificar"
$ws.Cells.Item(261, 9).Value = "Segunda"
$ws.Cells.Item(261, 10).Value = 300
$ws.Cells.Item(261, 11).Value = 600
$ws.Cells.Item(261, 12).Value = 600
$ws.Cells.Item(261, 13).Value = 600
$ws.Cells.Item(261, 14).Value = "`$/unidad"
$ws.Cells.Item(261, 15).Value = "Región del Maule"
$ws.Cells.Item(261, 16).Value = 600
$ws.Cells.Item(261, 17).Value = 1
$ws.Cells.Item(261, 18).Value = "Hortaliza"
